$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# --- Update the VOTE / COMMENT related API rows -------------------------
# Row 13: /post/:title/upvote -> /pages/upvote/:title
$ws.Range("C13").Value = "/pages/upvote/:title"

# Row 14: /post/:title/downvote -> /pages/downvote/:title
$ws.Range("C14").Value = "/pages/downvote/:title"

# Row 15: /post/:title/comment -> /pages/comment/:title
$ws.Range("C15").Value = "/pages/comment/:title"

# Row 16: /post/:post_title/:comment/upvote -> /pages/:title/comments/upvote
$ws.Range("E16").Value = "title"
$ws.Range("F16").Value = "comment_id"
$ws.Range("C16").Value = "/pages/:title/comments/upvote"

# Row 17: /post/:post_title/:comment/downvote -> /pages/:title/comments/downvote
$ws.Range("C17").Value = "/pages/:title/comments/downvote"
$ws.Range("E17").Value = "title"
$ws.Range("F17").Value = "comment_id"

# Row 18: /post/:post_title/:comment/reply -> /pages/:title/comments/reply
$ws.Range("C18").Value = "/pages/:title/comments/reply"
$ws.Range("E18").Value = "title"

# Row 19: /post/:post_title/:comment -> /pages/:title/comments
$ws.Range("C19").Value = "/pages/:title/comments"
$ws.Range("E19").Value = "title"
# F19 is a brand new cell on this row; copy the neighbouring cell's format
# (style index) before filling in its value.
$ws.Range("G19").Copy() | Out-Null
$ws.Range("F19").PasteSpecial(-4122) | Out-Null  # xlPasteFormats
$ws.Range("F19").Value = "comment_id"

# --- Row heights were recomputed by Excel across the whole table --------
for ($r = 1; $r -le 19; $r++) {
    $ws.Rows.Item($r).RowHeight = 13.2
}

# --- Move the active selection to C16 ------------------------------------
$ws.Range("C16").Select() | Out-Null
